$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Match the text number format used on the date columns before assigning
# the values, so the date-like strings are stored as text, not dates.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("K3").NumberFormat = "@"

# Add a new data row (row 3) mirroring the structure of row 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Gare"
$ws.Range("C3").Value = "St-Malo"
$ws.Range("D3").Value = "À partir de"
$ws.Range("E3").Value = "01/12/2017"
$ws.Range("H3").Value = "gare"
$ws.Range("I3").Value = "Rennes"
$ws.Range("J3").Value = "À partir de"
$ws.Range("K3").Value = "05/12/2017"
$ws.Range("F3").Value = "12h"
$ws.Range("L3").Value = "13h"

# Update selection to the newly added cell
$ws.Range("K3").Select()

# Configure page setup as applied by Excel
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
